$wb = $excel.ActiveWorkbook

# --- Add "same_elements" sheet (after the last existing sheet) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sameElements = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$sameElements.Name = "same_elements"

$headers = @(5, 50, 500, 5000, 50000, 500000)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $sameElements.Cells.Item(1, $i + 2).Value = $headers[$i]
}

$labels = @("byte", "int", "string", "date")
for ($i = 0; $i -lt $labels.Length; $i++) {
    $sameElements.Cells.Item($i + 2, 1).Value = $labels[$i]
}

$sameElements.Range("B2").Value = 0
$sameElements.Range("C2").Value = 0
$sameElements.Range("D2").Value = 0
$sameElements.Range("E2").Value = 0.001
$sameElements.Range("F2").Value = 0
$sameElements.Range("G2").Value = 0

$sameElements.Range("B3").Value = 0
$sameElements.Range("C3").Value = 0.001
$sameElements.Range("D3").Value = 0
$sameElements.Range("E3").Value = 0.001
$sameElements.Range("F3").Value = 0
$sameElements.Range("G3").Value = 0

$sameElements.Range("B4").Value = 0
$sameElements.Range("C4").Value = 0
$sameElements.Range("D4").Value = 0
$sameElements.Range("E4").Value = 1.0009999999999999 / 1000
$sameElements.Range("F4").Value = 0
$sameElements.Range("G4").Value = 0

$sameElements.Range("B5").Value = 0
$sameElements.Range("C5").Value = 0
$sameElements.Range("D5").Value = 0
$sameElements.Range("E5").Value = 1.0009999999999999 / 1000
$sameElements.Range("F5").Value = 0
$sameElements.Range("G5").Value = 0

# --- Add "partly_same" sheet (after "same_elements") ---
$partlySame = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sameElements)
$partlySame.Name = "partly_same"

for ($i = 0; $i -lt $headers.Length; $i++) {
    $partlySame.Cells.Item(1, $i + 2).Value = $headers[$i]
}
for ($i = 0; $i -lt $labels.Length; $i++) {
    $partlySame.Cells.Item($i + 2, 1).Value = $labels[$i]
}

$partlySame.Range("B2").Value = 0
$partlySame.Range("C2").Value = 0
$partlySame.Range("D2").Value = 0.013
$partlySame.Range("E2").Value = 1.89917
$partlySame.Range("F2").Value = 0
$partlySame.Range("G2").Value = 0

$partlySame.Range("B3").Value = 0
$partlySame.Range("C3").Value = 0
$partlySame.Range("D3").Value = 0.014002
$partlySame.Range("E3").Value = 1.94217
$partlySame.Range("F3").Value = 0
$partlySame.Range("G3").Value = 0

$partlySame.Range("B4").Value = 0
$partlySame.Range("C4").Value = 0
$partlySame.Range("D4").Value = 0.018001
$partlySame.Range("E4").Value = 2.2622
$partlySame.Range("F4").Value = 0
$partlySame.Range("G4").Value = 0

$partlySame.Range("B5").Value = 0
$partlySame.Range("C5").Value = 0
$partlySame.Range("D5").Value = 0.017002
$partlySame.Range("E5").Value = 1.911169
$partlySame.Range("F5").Value = 0
$partlySame.Range("G5").Value = 0

# partly_same becomes the active tab with a selection on K13
$partlySame.Activate()
$partlySame.Range("K13").Select()
